$d = $word.ActiveDocument

# Grab the full WordprocessingML for the document body as text so we can
# perform the structural edit (removing proofing-error bookmarks and
# adding the two new notes at the end) at the XML level, then hand the
# whole thing back to Word so it can re-parse/re-render it (which also
# naturally merges runs that used to be split only to host a <w:proofErr/>
# marker).
$xml = $d.Content.WordOpenXML

# 1) Drop every spell/grammar proofing-error marker. These are inert,
#    self-closing elements (<w:proofErr w:type="spellStart"/> etc.) that
#    don't carry any text, so simply deleting them is safe.
$xml = $xml -replace '<w:proofErr [^>]*/>', ''

# 2) The final paragraph of the document only contains the _GoBack
#    bookmark. Replace it with three paragraphs: the existing bookmark
#    paragraph becomes two new bullet points ("Zu beachten..." and
#    "Manchmal..."), with the bookmark now living in the second one,
#    followed by a trailing space run (matching how Word represents the
#    insertion point there).
$newTail = '<w:p w:rsidR="00186DC0" w:rsidRPr="00E83A76" w:rsidRDefault="00186DC0" w:rsidP="00450B04"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Zu beachten: ganze viele Nullpointer Checks! Man will in Tests nicht immer ALLE Bestandteile des zu testenden Objekts zusammensetzen m' + [char]0xFC + 'ssen, um nur einen Teil davon zu testen!</w:t></w:r></w:p><w:p w:rsidR="00186DC0" w:rsidRPr="00E83A76" w:rsidRDefault="00186DC0" w:rsidP="00450B04"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Manchmal zerbricht man sich den Kopf, wieso ein Test nicht l' + [char]0xE4 + 'uft. L' + [char]0xF6 + 'sung: einfach ein klein wenig l' + [char]0xE4 + 'nger warten, manchmal brauchen Dinge l' + [char]0xE4 + 'nger als man erwartet</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'

$tailPattern = '<w:p [^>]*w:rsidR="00186DC0"[^>]*><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

if (-not [regex]::IsMatch($xml, $tailPattern)) {
    throw "Could not locate the trailing bookmark paragraph to replace"
}

$xml = $xml -replace $tailPattern, $newTail

# Push the rebuilt OOXML back into the document, replacing the whole body.
$d.Content.InsertXML($xml)
